$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 100
$ws.Range("A100").Value = 99
$ws.Range("B100").Value = 12
$ws.Range("C100").Value = 'Dancing Queen'
$ws.Range("D100").Value = 'Gain 12% towards your damage for each level - for a total of +120% damage.'
$ws.Range("E100").Value = 1
$ws.Range("I100").Value = 0.012

# Row 101
$ws.Range("A101").Value = 100
$ws.Range("B101").Value = 12
$ws.Range("C101").Value = 'Confusing Day Dream Dance'
$ws.Range("D101").Value = 'Over time, you will reduce the enemies abilities to heal, do damage and so on, listed under Reductions Section, by 100%. You will also gain 100%, over time, in your base damage stat.'
$ws.Range("E101").Value = 12
$ws.Range("N101").Value = 0.01
$ws.Range("P101").Value = 0.01
$ws.Range("Q101").Value = 0.01
$ws.Range("R101").Value = 0.01
$ws.Range("S101").Value = 0.01
$ws.Range("T101").Value = 0.01

# Row 102
$ws.Range("A102").Value = 101
$ws.Range("B102").Value = 12
$ws.Range("C102").Value = 'Self Loving Dance'
$ws.Range("D102").Value = 'Increase your healing to 200% over time. You will also gain 140% over time, towards your health and 50% towards your damage stat'
$ws.Range("E102").Value = 24
$ws.Range("K102").Value = 0.02
$ws.Range("M102").Value = 0.014
$ws.Range("N102").Value = 0.005

# Row 103
$ws.Range("A103").Value = 102
$ws.Range("B103").Value = 12
$ws.Range("C103").Value = 'Deadly strip tease'
$ws.Range("D103").Value = 'You will deal 5,000 +25% of your damage stat growing by 50 points per level for a total of 10,000 Damage +25% of your base damage stat. You must use ATTACK for this to trigger.'
$ws.Range("E103").Value = 36
$ws.Range("F103").Value = 5000
$ws.Range("G103").Value = 50
$ws.Range("H103").Value = 0.25
$ws.Range("O103").Value = 'attack'

# Row 104
$ws.Range("A104").Value = 103
$ws.Range("B104").Value = 12
$ws.Range("C104").Value = 'Spell Weavers Dance'
$ws.Range("D104").Value = 'As you weave your magic into your dance you will do 10,000 + 18% of your damage stat as damage, growing by 100 over time for an additional 10,000 damage. Your spell damage will grow by +75% over time, while the enemies spell evasion will fall by 50% over time. You must use CAST attack for this to trigger.'
$ws.Range("E104").Value = 48
$ws.Range("F104").Value = 10000
$ws.Range("G104").Value = 100
$ws.Range("H104").Value = 0.18
$ws.Range("L104").Value = 0.0075
$ws.Range("O104").Value = 'cast'
$ws.Range("P104").Value = 0.005

# Row 105
$ws.Range("A105").Value = 104
$ws.Range("B105").Value = 12
$ws.Range("C105").Value = 'Double Dipping Tips'
$ws.Range("D105").Value = 'As you dance to steal your enemies tips, you will do 15,000 damage + 12% of your damage stat while growing this damage by 300 over time, for a total of 30,000 extra damage. You will also raise your damage by 200%, your spell damage by 100% and your damage stat by 50%, over time, Damage is only done when you use ATTACK AND CAST'
$ws.Range("E105").Value = 60
$ws.Range("F105").Value = 15000
$ws.Range("G105").Value = 300
$ws.Range("H105").Value = 0.12
$ws.Range("I105").Value = 0.02
$ws.Range("L105").Value = 0.01
$ws.Range("N105").Value = 0.005
$ws.Range("O105").Value = 'attack_and_cast'

# Row 106
$ws.Range("A106").Value = 105
$ws.Range("B106").Value = 12
$ws.Range("C106").Value = 'Rage Induced Dance'
$ws.Range("D106").Value = 'Reduces the affix damage and the resistances and enemy has by 100% over time. Deal 30,000 Damage growing by an additional 30,000 damage with +50% of your damage stat while you grow your damage stat by 200% over time. Damage is only dealt when using ATTACK.'
$ws.Range("E106").Value = 70
$ws.Range("F106").Value = 30000
$ws.Range("G106").Value = 300
$ws.Range("H106").Value = 0.5
$ws.Range("N106").Value = 0.02
$ws.Range("O106").Value = 'attack'
$ws.Range("Q106").Value = 0.01
$ws.Range("T106").Value = 0.01

# Row 107
$ws.Range("A107").Value = 106
$ws.Range("B107").Value = 12
$ws.Range("C107").Value = 'The Churches Dance'
$ws.Range("D107").Value = 'Reduce the enemies healing by 100% over time, While increasing your own healing and health by 200%. Deal damage of 50,000 increasing by 500 for an additional 50,000 damage while applying 25% of your damage stat towards the over all damage. Damage only procs if you use CAST AND ATTACK.'
$ws.Range("E107").Value = 80
$ws.Range("F107").Value = 50000
$ws.Range("G107").Value = 500
$ws.Range("H107").Value = 0.25
$ws.Range("K107").Value = 0.02
$ws.Range("M107").Value = 0.02
$ws.Range("O107").Value = 'cast_and_attack'
$ws.Range("R107").Value = 0.01

# Row 108
$ws.Range("A108").Value = 107
$ws.Range("B108").Value = 12
$ws.Range("C108").Value = 'Wedding Dance'
$ws.Range("D108").Value = 'Reduce all listed reductions by 100%, while dealing 100,000 Damage growing by an additional 100,000 damage over time and applying 20% of your damage stat to the damage while using ANY attack. Grow your damage by 250% over time and your base damage stat by 250% over time.'
$ws.Range("E108").Value = 90
$ws.Range("F108").Value = 100000
$ws.Range("G108").Value = 1000
$ws.Range("H108").Value = 0.2
$ws.Range("I108").Value = 0.25
$ws.Range("N108").Value = 0.025
$ws.Range("O108").Value = 'any'
$ws.Range("P108").Value = 0.01
$ws.Range("Q108").Value = 0.01
$ws.Range("R108").Value = 0.01
$ws.Range("S108").Value = 0.01
$ws.Range("T108").Value = 0.01

# Column width adjustments (approximate best-fit widths after new data added)
$ws.Columns.Item(1).ColumnWidth = 4.57
$ws.Columns.Item(2).ColumnWidth = 16.425
$ws.Columns.Item(3).ColumnWidth = 31.707
$ws.Columns.Item(4).ColumnWidth = 388.048
$ws.Columns.Item(5).ColumnWidth = 30.564
$ws.Columns.Item(6).ColumnWidth = 19.995
$ws.Columns.Item(7).ColumnWidth = 42.418
$ws.Columns.Item(8).ColumnWidth = 48.274
$ws.Columns.Item(9).ColumnWidth = 18.71
$ws.Columns.Item(10).ColumnWidth = 13.997
$ws.Columns.Item(11).ColumnWidth = 19.995
$ws.Columns.Item(12).ColumnWidth = 25.851
$ws.Columns.Item(13).ColumnWidth = 12.854
$ws.Columns.Item(14).ColumnWidth = 30.564
$ws.Columns.Item(15).ColumnWidth = 24.708
$ws.Columns.Item(16).ColumnWidth = 16.425
$ws.Columns.Item(17).ColumnWidth = 26.993
$ws.Columns.Item(18).ColumnWidth = 21.138
$ws.Columns.Item(19).ColumnWidth = 18.71
$ws.Columns.Item(20).ColumnWidth = 24.708
